# Weekly data refresh: a new price record (week of 2023-08-28) is inserted
# at row 462, pushing every subsequent record down by one row. The oldest
# record (previously the last row, 506) ends up duplicated onto the new
# last row, 507.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 462; this shifts rows 462-506 down to 463-507
# (and carries the row's number formatting, e.g. the date style in column D).
$ws.Rows.Item(462).Insert()

# Populate the newly inserted row 462 with the new weekly record.
$ws.Cells.Item(462, 1).Value = 4
$ws.Cells.Item(462, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(462, 3).Value = "Los Lagos"
$ws.Cells.Item(462, 4).Value = 45166
$ws.Cells.Item(462, 5).Value = 10
$ws.Cells.Item(462, 6).Value = 100112037
$ws.Cells.Item(462, 7).Value = "Cebollín"
$ws.Cells.Item(462, 8).Value = "Sin especificar"
$ws.Cells.Item(462, 9).Value = "Primera"
$ws.Cells.Item(462, 10).Value = 70
$ws.Cells.Item(462, 11).Value = 6000
$ws.Cells.Item(462, 12).Value = 6000
$ws.Cells.Item(462, 13).Value = 6000
$ws.Cells.Item(462, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(462, 15).Value = "Región Metropolitana"
$ws.Cells.Item(462, 16).Value = 167
$ws.Cells.Item(462, 17).Value = 36
$ws.Cells.Item(462, 18).Value = "Hortaliza"
